# Auto-generated edit script: updates market-price derived columns
# (H/I/J/K/L/M/N) across 8 sheets to match refreshed scrape data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 301.25
$ws.Range("I8").Value = 304.57144
$ws.Range("K8").Value = 913.71432
$ws.Range("M8").Value = -774.71432
$ws.Range("H40").Value = 1814.6154
$ws.Range("I40").Value = 1236
$ws.Range("K40").Value = 1236
$ws.Range("M40").Value = -1061
$ws.Range("H74").Value = 3703.818
$ws.Range("I74").Value = 3703.818
$ws.Range("K74").Value = 3703.818
$ws.Range("M74").Value = -2767.818
$ws.Range("H77").Value = 3703.818
$ws.Range("I77").Value = 3703.818
$ws.Range("K77").Value = 18519.09
$ws.Range("M77").Value = -13839.09
$ws.Range("H100").Value = 3195.7273
$ws.Range("I100").Value = 3139.3333
$ws.Range("K100").Value = 3139.3333
$ws.Range("M100").Value = -2598.3333
$ws.Range("H125").Value = 1097.4
$ws.Range("J125").Value = 1097.4
$ws.Range("L125").Value = 9876.6
$ws.Range("N125").Value = -14796.6
$ws.Range("I137").Value = 2052.1667
$ws.Range("J137").Value = 2740.8572
$ws.Range("K137").Value = 6156.500100000001
$ws.Range("L137").Value = 8222.571599999999
$ws.Range("M137").Value = -3606.500100000001
$ws.Range("N137").Value = -13322.5716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4917
$ws.Range("I32").Value = 3827.2173
$ws.Range("K32").Value = 3827.2173
$ws.Range("M32").Value = -3540.2173
$ws.Range("H45").Value = 2255
$ws.Range("I45").Value = 2132.5
$ws.Range("K45").Value = 2132.5
$ws.Range("M45").Value = -1755.5
$ws.Range("H74").Value = 1636.7059
$ws.Range("I74").Value = 1521.9333
$ws.Range("J74").Value = 2497.5
$ws.Range("K74").Value = 1521.9333
$ws.Range("L74").Value = 2497.5
$ws.Range("M74").Value = -647.9332999999999
$ws.Range("N74").Value = -4245.5
$ws.Range("H77").Value = 1636.7059
$ws.Range("I77").Value = 1521.9333
$ws.Range("J77").Value = 2497.5
$ws.Range("K77").Value = 7609.666499999999
$ws.Range("L77").Value = 12487.5
$ws.Range("M77").Value = -3241.666499999999
$ws.Range("N77").Value = -21223.5
$ws.Range("H102").Value = 1864.3334
$ws.Range("I102").Value = 1864.3334
$ws.Range("K102").Value = 1864.3334
$ws.Range("M102").Value = -242.3334
$ws.Range("H131").Value = 65356.25
$ws.Range("J131").Value = 65356.25
$ws.Range("L131").Value = 65356.25
$ws.Range("N131").Value = -75436.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3736.625
$ws.Range("I20").Value = 2919
$ws.Range("K20").Value = 2919
$ws.Range("M20").Value = -2672
$ws.Range("H134").Value = 5275.8887
$ws.Range("I134").Value = 5532.357
$ws.Range("K134").Value = 16597.071
$ws.Range("M134").Value = -14062.071

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 7894.8
$ws.Range("I86").Value = 7491.3335
$ws.Range("J86").Value = 8500
$ws.Range("K86").Value = 7491.3335
$ws.Range("L86").Value = 8500
$ws.Range("M86").Value = -6368.3335
$ws.Range("N86").Value = -10746
$ws.Range("H89").Value = 7894.8
$ws.Range("I89").Value = 7491.3335
$ws.Range("J89").Value = 8500
$ws.Range("K89").Value = 37456.6675
$ws.Range("L89").Value = 42500
$ws.Range("M89").Value = -31840.6675
$ws.Range("N89").Value = -53732
$ws.Range("H134").Value = 3632.9167
$ws.Range("I134").Value = 3632.9167
$ws.Range("K134").Value = 10898.7501
$ws.Range("M134").Value = -8363.750100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1470471.2
$ws.Range("I4").Value = 56993.133
$ws.Range("K4").Value = 170979.399
$ws.Range("M4").Value = -170867.399
$ws.Range("H6").Value = 378.66666
$ws.Range("I6").Value = 378.66666
$ws.Range("K6").Value = 1135.99998
$ws.Range("M6").Value = -1022.99998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 39000
$ws.Range("J27").Value = 39000
$ws.Range("L27").Value = 39000
$ws.Range("N27").Value = -39332
$ws.Range("H80").Value = 2631.077
$ws.Range("I80").Value = 1975.4
$ws.Range("K80").Value = 1975.4
$ws.Range("M80").Value = -977.4000000000001
$ws.Range("H83").Value = 2631.077
$ws.Range("I83").Value = 1975.4
$ws.Range("K83").Value = 9877
$ws.Range("M83").Value = -4885
$ws.Range("H122").Value = 3783.6667
$ws.Range("I122").Value = 2540.4
$ws.Range("K122").Value = 7621.200000000001
$ws.Range("M122").Value = -5171.200000000001
$ws.Range("H132").Value = 5374.3335
$ws.Range("I132").Value = 4494.4443
$ws.Range("J132").Value = 8014
$ws.Range("K132").Value = 13483.3329
$ws.Range("L132").Value = 24042
$ws.Range("M132").Value = -10953.3329
$ws.Range("N132").Value = -29102

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H21").Value = 9998.5
$ws.Range("I21").Value = 9997
$ws.Range("J21").Value = 10000
$ws.Range("K21").Value = 9997
$ws.Range("L21").Value = 10000
$ws.Range("N21").Value = -10348
$ws.Range("M21").Value = -9823
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("H40").Value = 3470
$ws.Range("I40").Value = 3470
$ws.Range("K40").Value = 3470
$ws.Range("M40").Value = -3334
$ws.Range("H136").Value = 29413564
$ws.Range("I136").Value = 1562.6428
$ws.Range("K136").Value = 4687.928400000001
$ws.Range("M136").Value = -2137.928400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 15000
$ws.Range("I12").Value = 15000
$ws.Range("K12").Value = 15000
$ws.Range("M12").Value = -14858
$ws.Range("H29").Value = 70000
$ws.Range("I29").Value = 50000
$ws.Range("J29").Value = 90000
$ws.Range("K29").Value = 50000
$ws.Range("L29").Value = 90000
$ws.Range("M29").Value = -49710
$ws.Range("N29").Value = -90580
$ws.Range("H44").Value = 37
$ws.Range("I44").Value = 37
$ws.Range("K44").Value = 37
$ws.Range("M44").Value = 517
$ws.Range("H81").Value = 4316.6
$ws.Range("I81").Value = 2266.5
$ws.Range("K81").Value = 4533
$ws.Range("M81").Value = -3472
$ws.Range("H84").Value = 4316.6
$ws.Range("I84").Value = 2266.5
$ws.Range("K84").Value = 22665
$ws.Range("M84").Value = -17361
$ws.Range("H126").Value = 3581.3635
$ws.Range("I126").Value = 3711.625
$ws.Range("K126").Value = 11134.875
$ws.Range("M126").Value = -8664.875
$ws.Range("H130").Value = 40497.5
$ws.Range("J130").Value = 40497.5
$ws.Range("L130").Value = 40497.5
$ws.Range("N130").Value = -50537.5
$ws.Range("H132").Value = 6136.15
$ws.Range("I132").Value = 3837.3572
$ws.Range("K132").Value = 11512.0716
$ws.Range("M132").Value = -8982.071599999999
$ws.Range("H136").Value = 6781.2964
$ws.Range("I136").Value = 2819.5
$ws.Range("K136").Value = 8458.5
$ws.Range("M136").Value = -5908.5

